$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Rows 21/22: BitcoinCash and Uniswap swapped places with new values
$ws.Cells.Item(21, 2).Value = "BitcoinCash"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue 21 4 "496.45"
$ws.Cells.Item(21, 5).Value = "  +3.85%  "

$ws.Cells.Item(22, 2).Value = "Uniswap"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue 22 4 "10.72"
$ws.Cells.Item(22, 5).Value = "  -3.98%  "

# Remaining Price (D) / Volume(1h) (E) updates
Set-TextValue 2 4 "71.263.82"
$ws.Cells.Item(2, 5).Value = "  +0.82%  "
Set-TextValue 3 4 "3.842.19"
$ws.Cells.Item(3, 5).Value = "  +1.19%  "
Set-TextValue 4 4 "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
Set-TextValue 5 4 "715.19"
$ws.Cells.Item(5, 5).Value = "  +2.53%  "
Set-TextValue 6 4 "172.80"
$ws.Cells.Item(6, 5).Value = "  +0.21%  "
Set-TextValue 7 4 "3.839.40"
$ws.Cells.Item(7, 5).Value = "  +1.25%  "
Set-TextValue 8 4 "1.00"
$ws.Cells.Item(8, 5).Value = "  -0.02%  "
Set-TextValue 9 4 "0.526"
$ws.Cells.Item(9, 5).Value = "  -0.07%  "
Set-TextValue 10 4 "0.163"
$ws.Cells.Item(10, 5).Value = "  +0.60%  "
Set-TextValue 11 4 "7.33"
$ws.Cells.Item(11, 5).Value = "  +0.99%  "
$ws.Cells.Item(12, 5).Value = "  +0.66%  "
Set-TextValue 13 4 "0.0000257"
$ws.Cells.Item(13, 5).Value = "  -0.34%  "
Set-TextValue 14 4 "36.77"
$ws.Cells.Item(14, 5).Value = "  +1.75%  "
Set-TextValue 15 4 "4.498.26"
$ws.Cells.Item(15, 5).Value = "  +1.37%  "
Set-TextValue 16 4 "3.847.23"
$ws.Cells.Item(16, 5).Value = "  +1.61%  "
Set-TextValue 17 4 "71.273.42"
$ws.Cells.Item(17, 5).Value = "  +0.79%  "
Set-TextValue 18 4 "7.23"
$ws.Cells.Item(18, 5).Value = "  +0.86%  "
$ws.Cells.Item(19, 5).Value = "  +0.74%  "
Set-TextValue 20 4 "17.41"
$ws.Cells.Item(20, 5).Value = "  -1.78%  "
Set-TextValue 23 4 "0.726"
$ws.Cells.Item(23, 5).Value = "  +2.29%  "
Set-TextValue 24 4 "85.26"
$ws.Cells.Item(24, 5).Value = "  +1.79%  "
$ws.Cells.Item(25, 5).Value = "  +3.28%  "
Set-TextValue 26 4 "10.64"
$ws.Cells.Item(26, 5).Value = "  +2.46%  "
$ws.Cells.Item(27, 5).Value = "  -1.02%  "
$ws.Cells.Item(28, 5).Value = "  +3.85%  "
Set-TextValue 29 4 "2.10"
$ws.Cells.Item(29, 5).Value = "  -2.28%  "
$ws.Cells.Item(30, 5).Value = "  +0.05%  "
$ws.Cells.Item(31, 5).Value = "  +0.28%  "
$ws.Cells.Item(32, 5).Value = "  -1.61%  "
Set-TextValue 33 4 "29.48"
$ws.Cells.Item(33, 5).Value = "  +0.37%  "
$ws.Cells.Item(34, 5).Value = "  -6.60%  "
Set-TextValue 35 4 "9.21"
$ws.Cells.Item(35, 5).Value = "  -0.09%  "
Set-TextValue 36 4 "3.808.73"
$ws.Cells.Item(36, 5).Value = "  +1.54%  "
Set-TextValue 37 4 "0.998"
$ws.Cells.Item(37, 5).Value = "  -0.32%  "
Set-TextValue 38 4 "0.103"
$ws.Cells.Item(38, 5).Value = "  +0.61%  "
Set-TextValue 39 4 "6.03"
$ws.Cells.Item(39, 5).Value = "  +0.77%  "
$ws.Cells.Item(40, 5).Value = "  +6.03%  "
Set-TextValue 41 4 "3.36"
$ws.Cells.Item(41, 5).Value = "  -1.14%  "
Set-TextValue 42 4 "2.27"
$ws.Cells.Item(42, 5).Value = "  +1.03%  "
$ws.Cells.Item(44, 5).Value = "  +0.16%  "
Set-TextValue 45 4 "0.000319"
$ws.Cells.Item(45, 5).Value = "  -2.39%  "
Set-TextValue 46 4 "163.82"
$ws.Cells.Item(46, 5).Value = "  +0.05%  "
Set-TextValue 47 4 "48.71"
$ws.Cells.Item(47, 5).Value = "  -0.28%  "
Set-TextValue 48 4 "419.99"
$ws.Cells.Item(48, 5).Value = "  +2.78%  "
Set-TextValue 49 4 "1.40"
$ws.Cells.Item(49, 5).Value = "  +1.82%  "
Set-TextValue 50 4 "8.63"
$ws.Cells.Item(50, 5).Value = "  +0.88%  "
$ws.Cells.Item(51, 5).Value = "  -0.79%  "
